$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column width for column F (min=6 max=6) ---
# Target stored width is 12.5703125 characters; the engine quantizes
# ColumnWidth to 1/6-character steps, so 12.5 is the closest reachable value.
$ws.Columns.Item(6).ColumnWidth = 11.71

# --- Convert the N11:N17 and P11:P17 per-row formulas into shared formulas ---
$ws.Range("N11:N17").Formula = "=M11/A`$3"
$ws.Range("P11:P17").Formula = "=O11/B`$3"

# --- New shared strings used by the new "source gas" mixing table ---
# (written in this exact order so the shared-string table gets the same
# index assignment as the target: 55 source gas, 56 Calib ppm, 57 Air,
# 58 CO2, 59 Q total)
$ws.Range("E25").Value = "source gas"
$ws.Range("F25").Value = "Calib ppm"
$ws.Range("I25").Value = "Air"
$ws.Range("H25").Value = "CO2"
$ws.Range("G25").Value = "Q total"
$ws.Range("J25").Value = "Trace"

# --- Data + formulas for the new calibration carrier strength mixing table ---
$ws.Range("E26").Value = 1000
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 50

$ws.Range("E27").Value = 1000
$ws.Range("F27").Value = 250
$ws.Range("G27").Value = 50

$ws.Range("E28").Value = 1000
$ws.Range("F28").Value = 500
$ws.Range("G28").Value = 50

$ws.Range("E29").Value = 1000
$ws.Range("F29").Value = 750
$ws.Range("G29").Value = 50

# Row 26 formulas (kept as plain, non-shared formulas, matching source)
$ws.Range("H26").Formula = "=G26*F26/E26"
$ws.Range("I26").Formula = "=G26*(E26-F26)/E26"

# Rows 27:30 formulas written as one range-assignment each so the engine
# groups them into shared formulas (matches t="shared" ref="H27:H30" / "I27:I30").
# Row 30 itself has no data in the target, but the shared-formula's declared
# ref still spans through row 30 (a leftover from the original fill-down), so
# write the full range then clear the row-30 formula cells again.
$ws.Range("H27:H30").Formula = "=G27*F27/E27"
$ws.Range("I27:I30").Formula = "=G27*(E27-F27)/E27"
$ws.Range("H30").ClearContents()
$ws.Range("I30").ClearContents()

# --- Move the active selection to where the user finished editing ---
$ws.Range("J31").Select() | Out-Null
